# Apply weekly update to "Hortaliza, Agricola del Norte S.A. de Arica - Betarraga"
# A new weekly price record is inserted at rows 111-112 (date 2021-09-01 / serial 44455),
# pushing every later weekly record down by one pair of rows (113..164), and the
# previously last record (old rows 161-162) is appended as brand-new rows 163-164.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Fecha (D) and the price columns (J volumen, K/L/M precios, P precio/kg) ---
# for the rows whose values change as part of the downward shift / new-week insert.
$ws.Range("D111").Value = 44455
$ws.Range("J111").Value = 1200
$ws.Range("D112").Value = 44455
$ws.Range("J112").Value = 800
$ws.Range("D113").Value = 44414
$ws.Range("J113").Value = 700
$ws.Range("K113").Value = 500
$ws.Range("L113").Value = 550
$ws.Range("M113").Value = 525
$ws.Range("P113").Value = 131
$ws.Range("D114").Value = 44414
$ws.Range("J114").Value = 900
$ws.Range("K114").Value = 500
$ws.Range("L114").Value = 550
$ws.Range("M114").Value = 525
$ws.Range("P114").Value = 105
$ws.Range("D115").Value = 44389
$ws.Range("D116").Value = 44389
$ws.Range("J116").Value = 1000
$ws.Range("D117").Value = 44249
$ws.Range("J117").Value = 900
$ws.Range("K117").Value = 600
$ws.Range("L117").Value = 650
$ws.Range("M117").Value = 625
$ws.Range("P117").Value = 156
$ws.Range("D118").Value = 44249
$ws.Range("J118").Value = 1100
$ws.Range("K118").Value = 600
$ws.Range("L118").Value = 650
$ws.Range("M118").Value = 625
$ws.Range("P118").Value = 125
$ws.Range("D119").Value = 44270
$ws.Range("K119").Value = 650
$ws.Range("L119").Value = 700
$ws.Range("M119").Value = 675
$ws.Range("P119").Value = 169
$ws.Range("D120").Value = 44270
$ws.Range("J120").Value = 1200
$ws.Range("K120").Value = 650
$ws.Range("L120").Value = 700
$ws.Range("M120").Value = 675
$ws.Range("P120").Value = 135
$ws.Range("D121").Value = 44260
$ws.Range("J121").Value = 1000
$ws.Range("K121").Value = 600
$ws.Range("L121").Value = 650
$ws.Range("M121").Value = 625
$ws.Range("P121").Value = 156
$ws.Range("D122").Value = 44260
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 600
$ws.Range("L122").Value = 650
$ws.Range("M122").Value = 625
$ws.Range("P122").Value = 125
$ws.Range("D123").Value = 44312
$ws.Range("J123").Value = 900
$ws.Range("K123").Value = 750
$ws.Range("L123").Value = 850
$ws.Range("M123").Value = 800
$ws.Range("P123").Value = 200
$ws.Range("D124").Value = 44312
$ws.Range("J124").Value = 1000
$ws.Range("K124").Value = 750
$ws.Range("L124").Value = 850
$ws.Range("M124").Value = 800
$ws.Range("P124").Value = 160
$ws.Range("D125").Value = 44386
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 650
$ws.Range("L125").Value = 700
$ws.Range("M125").Value = 675
$ws.Range("P125").Value = 169
$ws.Range("D126").Value = 44386
$ws.Range("J126").Value = 1300
$ws.Range("K126").Value = 650
$ws.Range("L126").Value = 700
$ws.Range("M126").Value = 675
$ws.Range("P126").Value = 135
$ws.Range("D127").Value = 44264
$ws.Range("K127").Value = 600
$ws.Range("L127").Value = 650
$ws.Range("M127").Value = 625
$ws.Range("P127").Value = 156
$ws.Range("D128").Value = 44264
$ws.Range("J128").Value = 900
$ws.Range("K128").Value = 600
$ws.Range("L128").Value = 650
$ws.Range("M128").Value = 625
$ws.Range("P128").Value = 125
$ws.Range("D129").Value = 44301
$ws.Range("J129").Value = 800
$ws.Range("K129").Value = 750
$ws.Range("L129").Value = 800
$ws.Range("M129").Value = 775
$ws.Range("P129").Value = 194
$ws.Range("D130").Value = 44301
$ws.Range("J130").Value = 600
$ws.Range("K130").Value = 750
$ws.Range("L130").Value = 800
$ws.Range("M130").Value = 775
$ws.Range("P130").Value = 155
$ws.Range("D131").Value = 44433
$ws.Range("J131").Value = 600
$ws.Range("K131").Value = 550
$ws.Range("L131").Value = 600
$ws.Range("M131").Value = 575
$ws.Range("P131").Value = 144
$ws.Range("D132").Value = 44433
$ws.Range("J132").Value = 800
$ws.Range("K132").Value = 550
$ws.Range("L132").Value = 600
$ws.Range("M132").Value = 575
$ws.Range("P132").Value = 115
$ws.Range("D133").Value = 44221
$ws.Range("J133").Value = 1200
$ws.Range("L133").Value = 550
$ws.Range("M133").Value = 525
$ws.Range("P133").Value = 131
$ws.Range("D134").Value = 44221
$ws.Range("J134").Value = 1000
$ws.Range("L134").Value = 550
$ws.Range("M134").Value = 525
$ws.Range("P134").Value = 105
$ws.Range("D135").Value = 44397
$ws.Range("J135").Value = 900
$ws.Range("K135").Value = 500
$ws.Range("L135").Value = 600
$ws.Range("M135").Value = 550
$ws.Range("P135").Value = 138
$ws.Range("D136").Value = 44397
$ws.Range("J136").Value = 960
$ws.Range("K136").Value = 500
$ws.Range("L136").Value = 600
$ws.Range("M136").Value = 550
$ws.Range("P136").Value = 110
$ws.Range("D137").Value = 44273
$ws.Range("J137").Value = 500
$ws.Range("K137").Value = 600
$ws.Range("L137").Value = 650
$ws.Range("M137").Value = 625
$ws.Range("P137").Value = 156
$ws.Range("D138").Value = 44273
$ws.Range("J138").Value = 600
$ws.Range("K138").Value = 600
$ws.Range("L138").Value = 650
$ws.Range("M138").Value = 625
$ws.Range("P138").Value = 125
$ws.Range("D139").Value = 44438
$ws.Range("J139").Value = 1200
$ws.Range("K139").Value = 450
$ws.Range("L139").Value = 500
$ws.Range("M139").Value = 475
$ws.Range("P139").Value = 119
$ws.Range("D140").Value = 44438
$ws.Range("J140").Value = 1200
$ws.Range("K140").Value = 450
$ws.Range("L140").Value = 500
$ws.Range("M140").Value = 475
$ws.Range("P140").Value = 95
$ws.Range("D141").Value = 44286
$ws.Range("J141").Value = 800
$ws.Range("K141").Value = 750
$ws.Range("L141").Value = 800
$ws.Range("M141").Value = 775
$ws.Range("P141").Value = 194
$ws.Range("D142").Value = 44286
$ws.Range("J142").Value = 800
$ws.Range("K142").Value = 750
$ws.Range("L142").Value = 800
$ws.Range("M142").Value = 775
$ws.Range("P142").Value = 155
$ws.Range("D143").Value = 44351
$ws.Range("K143").Value = 550
$ws.Range("L143").Value = 600
$ws.Range("M143").Value = 575
$ws.Range("P143").Value = 144
$ws.Range("D144").Value = 44351
$ws.Range("K144").Value = 550
$ws.Range("L144").Value = 600
$ws.Range("M144").Value = 575
$ws.Range("P144").Value = 115
$ws.Range("D145").Value = 44365
$ws.Range("J145").Value = 900
$ws.Range("K145").Value = 600
$ws.Range("L145").Value = 650
$ws.Range("M145").Value = 625
$ws.Range("P145").Value = 156
$ws.Range("D146").Value = 44365
$ws.Range("J146").Value = 1000
$ws.Range("K146").Value = 600
$ws.Range("L146").Value = 650
$ws.Range("M146").Value = 625
$ws.Range("P146").Value = 125
$ws.Range("D147").Value = 44306
$ws.Range("J147").Value = 400
$ws.Range("K147").Value = 750
$ws.Range("L147").Value = 800
$ws.Range("M147").Value = 775
$ws.Range("P147").Value = 194
$ws.Range("D148").Value = 44306
$ws.Range("J148").Value = 400
$ws.Range("K148").Value = 750
$ws.Range("L148").Value = 800
$ws.Range("M148").Value = 775
$ws.Range("P148").Value = 155
$ws.Range("D149").Value = 44162
$ws.Range("J149").Value = 1600
$ws.Range("K149").Value = 350
$ws.Range("L149").Value = 400
$ws.Range("M149").Value = 375
$ws.Range("P149").Value = 94
$ws.Range("D150").Value = 44162
$ws.Range("J150").Value = 1450
$ws.Range("K150").Value = 350
$ws.Range("L150").Value = 400
$ws.Range("M150").Value = 375
$ws.Range("P150").Value = 75
$ws.Range("D151").Value = 44410
$ws.Range("K151").Value = 550
$ws.Range("L151").Value = 600
$ws.Range("M151").Value = 575
$ws.Range("P151").Value = 144
$ws.Range("D152").Value = 44410
$ws.Range("K152").Value = 550
$ws.Range("L152").Value = 600
$ws.Range("M152").Value = 575
$ws.Range("P152").Value = 115
$ws.Range("D153").Value = 44176
$ws.Range("J153").Value = 1200
$ws.Range("K153").Value = 350
$ws.Range("L153").Value = 400
$ws.Range("M153").Value = 375
$ws.Range("P153").Value = 94
$ws.Range("D154").Value = 44176
$ws.Range("J154").Value = 1200
$ws.Range("K154").Value = 350
$ws.Range("L154").Value = 400
$ws.Range("M154").Value = 375
$ws.Range("P154").Value = 75
$ws.Range("D155").Value = 44239
$ws.Range("J155").Value = 700
$ws.Range("M155").Value = 625
$ws.Range("P155").Value = 156
$ws.Range("D156").Value = 44239
$ws.Range("J156").Value = 1000
$ws.Range("M156").Value = 625
$ws.Range("P156").Value = 125
$ws.Range("D157").Value = 44376
$ws.Range("J157").Value = 750
$ws.Range("K157").Value = 600
$ws.Range("L157").Value = 650
$ws.Range("M157").Value = 630
$ws.Range("P157").Value = 158
$ws.Range("D158").Value = 44376
$ws.Range("J158").Value = 700
$ws.Range("K158").Value = 600
$ws.Range("L158").Value = 650
$ws.Range("M158").Value = 629
$ws.Range("P158").Value = 126
$ws.Range("D159").Value = 44358
$ws.Range("J159").Value = 1200
$ws.Range("K159").Value = 450
$ws.Range("L159").Value = 500
$ws.Range("M159").Value = 475
$ws.Range("P159").Value = 119
$ws.Range("D160").Value = 44358
$ws.Range("J160").Value = 1200
$ws.Range("K160").Value = 450
$ws.Range("L160").Value = 500
$ws.Range("M160").Value = 475
$ws.Range("P160").Value = 95
$ws.Range("D161").Value = 44425
$ws.Range("J161").Value = 700
$ws.Range("K161").Value = 500
$ws.Range("L161").Value = 550
$ws.Range("M161").Value = 525
$ws.Range("P161").Value = 131
$ws.Range("D162").Value = 44425
$ws.Range("J162").Value = 900
$ws.Range("K162").Value = 500
$ws.Range("L162").Value = 550
$ws.Range("M162").Value = 525
$ws.Range("P162").Value = 105

# --- Append the two brand-new rows 163 (Primera) and 164 (Segunda) at the end of the table ---
# These carry the data that used to live in rows 161-162 before the shift.

# Row 163
$ws.Range("A163").Value = 1
$ws.Range("B163").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C163").Value = "Arica y Parinacota"
$ws.Range("D163").Value = 44323
$ws.Range("D163").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E163").Value = 15
$ws.Range("F163").Value = 100114014
$ws.Range("G163").Value = "Betarraga"
$ws.Range("H163").Value = "Sin especificar"
$ws.Range("I163").Value = "Primera"
$ws.Range("J163").Value = 1200
$ws.Range("K163").Value = 700
$ws.Range("L163").Value = 750
$ws.Range("M163").Value = 725
$ws.Range("N163").Value = "$/paquete 4 unidades"
$ws.Range("O163").Value = "Región de Arica y Parinacota"
$ws.Range("P163").Value = 181
$ws.Range("Q163").Value = 4
$ws.Range("R163").Value = "Hortaliza"

# Row 164
$ws.Range("A164").Value = 1
$ws.Range("B164").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C164").Value = "Arica y Parinacota"
$ws.Range("D164").Value = 44323
$ws.Range("D164").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E164").Value = 15
$ws.Range("F164").Value = 100114014
$ws.Range("G164").Value = "Betarraga"
$ws.Range("H164").Value = "Sin especificar"
$ws.Range("I164").Value = "Segunda"
$ws.Range("J164").Value = 1000
$ws.Range("K164").Value = 700
$ws.Range("L164").Value = 750
$ws.Range("M164").Value = 725
$ws.Range("N164").Value = "$/paquete 5 unidades"
$ws.Range("O164").Value = "Región de Arica y Parinacota"
$ws.Range("P164").Value = 145
$ws.Range("Q164").Value = 5
$ws.Range("R164").Value = "Hortaliza"
